# Updated cryptos list values (Price / Volume(1h) columns) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.769.82"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").Value = "1.784.44"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "310.64"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.5107"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("D8").Value = "0.3859"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").Value = "0.07821"
$ws.Range("E9").Value = "  -8.06%  "
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").Value = "40.69"
$ws.Range("E11").Value = "  -2.79%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "6.197"
$ws.Range("D14").Value = "20.15"
$ws.Range("E14").Value = "  -4.25%  "
$ws.Range("D15").Value = "1.776.46"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").Value = "7.209"
$ws.Range("E16").Value = "  -4.34%  "
$ws.Range("D17").Value = "91.47"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("E18").Value = "  -5.74%  "
$ws.Range("D19").Value = "0.06561"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "1.002"
$ws.Range("E21").Value = "  -4.15%  "
$ws.Range("D22").Value = "5.908"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").Value = "27.823.42"
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.00"
$ws.Range("E24").Value = "  -4.19%  "
$ws.Range("D25").Value = "2.227"
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("D26").Value = "159.85"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -3.96%  "
$ws.Range("D28").Value = "1.984.21"
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("D29").Value = "2.371"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").Value = "123.27"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").Value = "0.1076"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").Value = "1.034"
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("D33").Value = "3.635"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "5.483"
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("D35").Value = "0.07067"
$ws.Range("E35").Value = "  -4.73%  "
$ws.Range("D36").Value = "0.02305"
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").Value = "8.792"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2120"
$ws.Range("E38").Value = "  -5.09%  "
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").Value = "4.998"
$ws.Range("E40").Value = "  -4.39%  "
$ws.Range("D41").Value = "0.6086"
$ws.Range("E41").Value = "  -3.78%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "1.153"
$ws.Range("E43").Value = "  -3.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.320"
$ws.Range("E44").Value = "  -5.67%  "
$ws.Range("D45").Value = "13.18"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("D46").Value = "0.5903"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").Value = "125.74"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("E50").Value = "  -4.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06850"
$ws.Range("E51").Value = "  -1.80%  "
